$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.833.56'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.55%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.345.75'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +8.51%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.84'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +7.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '620.90'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.52%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.18'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +6.66%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.383'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.63%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.03%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.342.13'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +8.57%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.797'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.06%  '

# Row 12
$ws.Range("E12").Value = '  +1.06%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '97.642.26'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +3.87%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.66'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.89%  '

# Row 15
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.973.33'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +8.92%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000245'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.79%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.46'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.09%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.345.77'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +9.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.59'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.35%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.71'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.93%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '478.61'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +8.48%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.85'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.31%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000206'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +8.64%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.08'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.82%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.66'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.51%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.49'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.35%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.91'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.14%  '

# Row 28
$ws.Range("E28").Value = '  +10.34%  '

# Row 29
$ws.Range("E29").Value = '  -0.22%  '

# Row 30
$ws.Range("E30").Value = '  +5.71%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.246'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.64%  '

# Row 32
$ws.Range("E32").Value = '  -0.71%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.20%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.11'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.87%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.11'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +6.65%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '519.53'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +8.85%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.150'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.86%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.23'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -6.28%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.93'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.46%  '

# Row 40
$ws.Range("E40").Value = '  +3.08%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.447'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.77%  '

# Row 42
$ws.Range("E42").Value = '  -0.75%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.63'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.88%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.789'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +17.16%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.21'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.53%  '

# Row 46
$ws.Range("E46").Value = '  +0.03%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.92'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.33%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.91'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +5.21%  '

# Row 49
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.48'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.33%  '

# Row 50
$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.36'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.99%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.48'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.23%  '
